# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 310
$ws1.Range("F4").Value  = 166
$ws1.Range("F5").Value  = 193
$ws1.Range("F6").Value  = 325
$ws1.Range("F8").Value  = 2161
$ws1.Range("F9").Value  = 372
$ws1.Range("F10").Value = 5299
$ws1.Range("F11").Value = 118
$ws1.Range("F12").Value = 355

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 310
$ws4.Range("F5").Value  = 166
$ws4.Range("F6").Value  = 193
$ws4.Range("F7").Value  = 325
$ws4.Range("F11").Value = 2161
$ws4.Range("F12").Value = 372
$ws4.Range("F13").Value = 5299
$ws4.Range("F14").Value = 118
$ws4.Range("F15").Value = 355
